# Fill marker info ("NAT") into column J for the rows that are
# currently missing it, then select A2:J45 with active cell A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows missing the "NAT" marker value in column J.
$rows = @(4, 5, 6, 7, 8, 9, 10, 22, 23, 26, 27, 28, 29, 30, 31, 32, 44, 45)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = "NAT"
}

# Update selection to match the committed state.
$ws.Range("A2:J45").Select()
